# Update 1.5.12: Changes to SIP, Import template MJD fix
#
# - Add a new "ТП" ("substation") lookup sheet right after "СИП",
#   holding the "Наименование ТП" header used as the list source.
# - Add a new "Питается от ТП" column (F) to the "СИП" sheet, with a
#   list-style data validation pointing at the new "ТП" sheet.

$wb = $excel.ActiveWorkbook

$wsSip = $wb.Worksheets.Item("СИП")

# --- 1. Extend "СИП" with the new "Питается от ТП" column ------------------
# Set this string first so it becomes the lower shared-string index (mirrors
# the authoring order captured in the diff).
$wsSip.Range("E1").Copy($wsSip.Range("F1"))
$wsSip.Range("F1").Value = "Питается от ТП"
$wsSip.Columns.Item(6).AutoFit()

# --- 2. Insert the new "ТП" sheet right after "СИП" -------------------------
$tp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSip)
$tp.Name = "ТП"
$tp.Range("A1").Value = "Наименование ТП"
$tp.Range("A1").Font.Bold = $true
$tp.Columns.Item(1).AutoFit()

# --- 3. List-validate the new column off of the "ТП" sheet -----------------
$wsSip.Range("F2:F1048576").Validation.Add(3, 1, 1, "=ТП!`$A`$2:`$A`$1048576")

# --- 4. Cosmetic: restore the builtin "Normal" style's localized name ------
$wb.Styles.Item(1).Name = "Обычный"

# --- 5. Restore the sheet view state captured in the diff -------------------
$wsSip.Select()
$excel.ActiveWindow.Zoom = 175
$wsSip.Range("B16").Select()
